$wb = $excel.ActiveWorkbook

# Sheet "展览" (exhibition) - first sheet in the original workbook order
$wsExhibition = $wb.Worksheets.Item("展览")
$wsExhibition.Range("F3").Value = 26
$wsExhibition.Range("F4").Value = 22
$wsExhibition.Range("F5").Value = 4126
$wsExhibition.Range("F7").Value = 55
$wsExhibition.Range("F8").Value = 259
$wsExhibition.Range("F9").Value = 26

# Sheet "全部类型" (all types) - fourth sheet in the original workbook order
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F7").Value = 26
$wsAll.Range("F8").Value = 22
$wsAll.Range("F9").Value = 4126
$wsAll.Range("F11").Value = 55
$wsAll.Range("F13").Value = 259
$wsAll.Range("F14").Value = 26
